$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.633908748626709
$ws.Range("B1").Value = 2.457040786743164
$ws.Range("C1").Value = 5.051864147186279
$ws.Range("D1").Value = 3.878090858459473
$ws.Range("E1").Value = 1.861559152603149
